# Update the res_line/loading_percent sheet with refreshed loading-percentage
# results for the "contingencies with rene fine" case (Case_5_0).
# Columns B-N for data rows 2-25 get new values; this also introduces
# previously-absent columns G and J (both now present with data: G is 0,
# J holds a constant contingency value).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @(13.15886764000083, 7.462048855499266, 2.693823771289425, 4.796789376553945, 45.0797500448222, 0, 4.815596535470627, 5.454530539833984, 4.251792669778333, 31.68432192086116, 5.612467426898421, 11.63687459312686, 6.389801838044705)
    3 = @(12.40043597336053, 6.958510051112619, 2.705803728429739, 4.633496378861523, 43.58542334654938, 0, 5.099677360783825, 5.691181207635335, 4.251792669778333, 30.80574928454257, 5.502814417243211, 11.04720000993261, 6.19539558278994)
    4 = @(11.91317053789679, 6.676214771095246, 2.711929350760895, 4.530311791993292, 42.63701189621509, 0, 5.279964119801652, 5.841773446557165, 4.251792669778333, 30.24755122618214, 5.432987916399219, 10.67351263452058, 6.074103073062481)
    5 = @(11.70844495900099, 6.565725660119027, 2.712655437323071, 4.48771121421294, 42.21082589191396, 0, 5.355600757926628, 5.90685742681891, 4.251792669778333, 29.99098292428068, 5.402486839867177, 10.51982425895091, 6.025967830592509)
    6 = @(11.67306129154653, 6.554461675668477, 2.710991317742037, 4.480812773594291, 42.10007230396383, 0, 5.369063089960862, 5.920726819980366, 4.251792669778333, 29.9180380285027, 5.395652166082106, 10.49577726031437, 6.020070195342524)
    7 = @(11.9075186252457, 6.693832907014096, 2.707123315689427, 4.530330650594821, 42.52451972238111, 0, 5.283173040238574, 5.850572611340237, 4.251792669778333, 30.16274719155235, 5.42788876444927, 10.67579672184053, 6.079185632251393)
    8 = @(12.89817325313328, 7.303966544330573, 2.691999558286466, 4.741837509992613, 44.43566893897093, 0, 4.9151636516664, 5.54510833506696, 4.251792669778333, 31.28145940263289, 5.569155175375623, 11.44132785459056, 6.330506954754989)
    9 = @(14.66983687978676, 8.534783404703008, 2.660546102403116, 5.132395941920403, 48.10267991981512, 0, 4.234816142891945, 4.973801813422744, 4.251792669778333, 33.45654326980743, 5.834451965012127, 12.8370884753058, 6.796532687917221)
    10 = @(15.85980336152158, 9.379258029418121, 2.612548163143793, 5.348213633107304, 50.11372137984916, 0, 3.790591474958236, 4.592347544680055, 4.251792669778333, 34.59239189280667, 5.971355015985972, 13.79899394686357, 7.060141554711806)
    11 = @(16.38205131255406, 9.747673319776794, 2.464652917968847, 5.012354403792022, 46.91554724064454, 0, 4.32980069957162, 4.523709711458487, 4.251792669778333, 32.24718235773513, 5.695782743644475, 14.23615681808133, 6.601350952752953)
    12 = @(16.58169032711004, 9.864219642154353, 2.380525649043913, 4.776975146118486, 43.93588626717033, 0, 5.340082709690272, 4.524733590070461, 4.251792669778333, 30.16191369453594, 5.513014520106521, 14.39979601285817, 6.162266768870196)
    13 = @(16.54810868128771, 9.818170843018773, 2.333852998214274, 4.613200563927247, 40.75551745336475, 0, 6.550225427034033, 4.592128514706866, 4.251792669778333, 28.02063844435328, 5.388220653471034, 14.37272170064251, 5.717834786921058)
    14 = @(16.41830124607314, 9.716300571509603, 2.321226286867731, 4.550003129239076, 38.40538671295267, 0, 7.463365892642869, 4.671112015101583, 4.251792669778333, 26.48449645047081, 5.339285064664357, 14.26851016127788, 5.410017546727598)
    15 = @(16.33573206082777, 9.661939467494554, 2.323491608367373, 4.536242647215169, 37.75544752378092, 0, 7.686319463990333, 4.708864888482359, 4.251792669778333, 26.07185336431888, 5.329350471815518, 14.20379214086595, 5.334137854911566)
    16 = @(15.85415223608508, 9.331318763833989, 2.34971759152325, 4.465912770413847, 37.36830542906188, 0, 7.57770955588499, 4.862373615341233, 4.251792669778333, 25.92209171386809, 5.289779646346656, 13.81563622009909, 5.303800329626883)
    17 = @(15.54889932099007, 9.129993960702725, 2.374873800264222, 4.434539443209374, 38.33869885498497, 0, 6.994652213316406, 4.939226753814424, 4.251792669778333, 26.62403291831062, 5.277539803673451, 13.57023339789343, 5.443545815115701)
    18 = @(15.3676464277229, 9.008089069257821, 2.4145133421013, 4.494906245552678, 40.69250099793052, 0, 6.003473836076186, 4.947029479861273, 4.251792669778333, 28.23314333787769, 5.331811676603123, 13.41928682584766, 5.756473292104371)
    19 = @(15.29569538697847, 8.987548634598159, 2.481362793325466, 4.709797824584957, 43.84515037294788, 0, 4.92787171821319, 4.908874077668508, 4.251792669778333, 30.37480935847638, 5.492724506245604, 13.36166931060378, 6.210587373256598)
    20 = @(15.54700756250463, 9.204702891719503, 2.612077180740929, 5.290982868895878, 49.28909279460176, 0, 3.912228188992446, 4.716891901894972, 4.251792669778333, 34.063152118269, 5.921080330071044, 13.56132918457266, 7.004107141004232)
    21 = @(16.42359074618386, 9.824805902931772, 2.60047182060187, 5.541226557260266, 51.55837176008252, 0, 3.532905612548928, 4.409809997545246, 4.251792669778333, 35.45078004309013, 6.093702736748138, 14.29405570050567, 7.305006320943579)
    22 = @(16.97965094185079, 10.20185919833921, 2.588838988200416, 5.672180377046955, 52.88041005988025, 0, 3.3007757190209, 4.207647736979992, 4.251792669778333, 36.25942621949985, 6.187235929342986, 14.81884749992137, 7.457400743488682)
    23 = @(16.6878041942406, 9.986540245310794, 2.600060158711839, 5.601961664337364, 52.28203285128195, 0, 3.422014678460747, 4.305230132990052, 4.251792669778333, 35.91126026512404, 6.142365965311726, 14.53840585272919, 7.370611263772894)
    24 = @(15.53584580718013, 9.171475946936166, 2.632128650472816, 5.33204105997512, 49.78536132959433, 0, 3.892577603575105, 4.698174342283017, 4.251792669778333, 34.42522019294025, 5.960811710753553, 13.54238028751481, 7.045109406816423)
    25 = @(14.20503725250564, 8.240501018898428, 2.662184487913013, 5.030163285880678, 46.96621963096387, 0, 4.417947855279126, 5.138207541993952, 4.251792669778333, 32.75109712874725, 5.756890199440372, 12.47719956452652, 6.68201112642396)
}

foreach ($rowNum in $newValues.Keys) {
    $rowValues = $newValues[$rowNum]
    $arr = New-Object "object[,]" 1,13
    for ($col = 0; $col -lt 13; $col++) {
        $arr[0, $col] = $rowValues[$col]
    }
    $startCell = $ws.Cells.Item($rowNum, 2)   # column B
    $endCell = $ws.Cells.Item($rowNum, 14)    # column N
    $rowRange = $ws.Range($startCell, $endCell)
    $rowRange.Value = $arr
}
